# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "69.558.72"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3
$ws.Range("D3").Value = "2.493.68"
$ws.Range("E3").Value = "  -0.91%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "569.24"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6
Set-TextValue $ws.Range("D6") "164.77"
$ws.Range("E6").Value = "  -1.23%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("E8").Value = "  -1.38%  "

# Row 9
$ws.Range("D9").Value = "2.491.87"
$ws.Range("E9").Value = "  -0.80%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.158"
$ws.Range("E10").Value = "  -1.65%  "

# Row 11
$ws.Range("E11").Value = "  -0.36%  "

# Row 12
$ws.Range("E12").Value = "  +0.25%  "

# Row 13
Set-TextValue $ws.Range("D13") "4.90"
$ws.Range("E13").Value = "  -0.49%  "

# Row 14
$ws.Range("D14").Value = "2.950.56"
$ws.Range("E14").Value = "  -0.99%  "

# Row 15
$ws.Range("D15").Value = "69.435.39"
$ws.Range("E15").Value = "  -0.45%  "

# Row 16
$ws.Range("E16").Value = "  -0.58%  "

# Row 17
Set-TextValue $ws.Range("D17") "24.28"
$ws.Range("E17").Value = "  -2.82%  "

# Row 18
$ws.Range("D18").Value = "2.496.87"
$ws.Range("E18").Value = "  -0.65%  "

# Row 19
Set-TextValue $ws.Range("D19") "11.18"
$ws.Range("E19").Value = "  -1.94%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.39"
$ws.Range("E20").Value = "  -4.66%  "

# Row 21
Set-TextValue $ws.Range("D21") "347.03"
$ws.Range("E21").Value = "  -0.80%  "

# Row 22
$ws.Range("E22").Value = "  -1.19%  "

# Row 23
$ws.Range("E23").Value = "  -3.98%  "

# Row 24
$ws.Range("E24").Value = "  -0.05%  "

# Row 25
Set-TextValue $ws.Range("D25") "69.65"
$ws.Range("E25").Value = "  -1.32%  "

# Row 26
Set-TextValue $ws.Range("D26") "3.91"
$ws.Range("E26").Value = "  -2.29%  "

# Row 27
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.621.58"
$ws.Range("E27").Value = "  -1.47%  "

# Row 28
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D28") "8.64"
$ws.Range("E28").Value = "  -2.78%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.28%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0873"
$ws.Range("E30").Value = "  -2.99%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.66"
$ws.Range("E31").Value = "  -3.09%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D32") "1.19"
$ws.Range("E32").Value = "  -5.60%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D33") "438.11"
$ws.Range("E33").Value = "  -5.92%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  +0.17%  "

# Row 35
$ws.Range("E35").Value = "  -2.09%  "

# Row 36
Set-TextValue $ws.Range("D36") "156.52"
$ws.Range("E36").Value = "  +0.32%  "

# Row 37
$ws.Range("E37").Value = "  -2.94%  "

# Row 38
Set-TextValue $ws.Range("D38") "19.07"
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
Set-TextValue $ws.Range("D39") "18.18"
$ws.Range("E39").Value = "  -2.03%  "

# Row 40
$ws.Range("E40").Value = "  +0.01%  "

# Row 41
$ws.Range("E41").Value = "  -1.32%  "

# Row 42
Set-TextValue $ws.Range("D42") "4.59"
$ws.Range("E42").Value = "  -3.87%  "

# Row 43
$ws.Range("E43").Value = "  -2.04%  "

# Row 44
$ws.Range("E44").Value = "  -5.77%  "

# Row 45
$ws.Range("E45").Value = "  -7.70%  "

# Row 46
Set-TextValue $ws.Range("D46") "138.93"
$ws.Range("E46").Value = "  -2.45%  "

# Row 47
$ws.Range("E47").Value = "  -1.57%  "

# Row 48
$ws.Range("E48").Value = "  -3.22%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0723"
$ws.Range("E49").Value = "  -0.98%  "

# Row 50
$ws.Range("E50").Value = "  -0.57%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0923"
$ws.Range("E51").Value = "  -0.75%  "
